# Fruta / hortaliza, semanal
# Insert 6 new weekly rows of "Damasco" data right after the existing row 85,
# shifting all subsequent rows down by 6 (old row 86 becomes row 92, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows starting at row 86 (pushes old rows 86:151 down to 92:157)
$ws.Rows("86:91").Insert()

# Common fixed values shared by every Damasco record in this sheet
$A = 9
$B = "Vega Central Mapocho de Santiago"
$C = "Metropolitana"
$E = 13
$F = "Fruta"
$G = 100103
$H = "Frutos de hueso (carozo)"
$I = 100103003
$J = "Damasco"

# New data rows (row number -> field values)
$newRows = @(
    @{ Row = 86; D = 44923; K = "Dina"; L = "Especial"; M = 290; N = 18000; O = 18000; P = 18000; Q = "`$/caja 15 kilos granel"; R = "Provincia de Los Andes"; S = 1200; T = 15 },
    @{ Row = 87; D = 44923; K = "Dina"; L = "Primera";  M = 250; N = 15000; O = 15000; P = 15000; Q = "`$/caja 15 kilos granel"; R = "Provincia de Los Andes"; S = 1000; T = 15 },
    @{ Row = 88; D = 44923; K = "Dina"; L = "Segunda";  M = 220; N = 12000; O = 12000; P = 12000; Q = "`$/caja 15 kilos granel"; R = "Provincia de Los Andes"; S =  800; T = 15 },
    @{ Row = 89; D = 44923; K = "Dina"; L = "Especial"; M = 290; N = 18000; O = 18000; P = 18000; Q = "`$/caja 15 kilos granel"; R = "Provincia de Los Andes"; S = 1200; T = 15 },
    @{ Row = 90; D = 44923; K = "Dina"; L = "Primera";  M = 250; N = 15000; O = 15000; P = 15000; Q = "`$/caja 15 kilos granel"; R = "Provincia de Los Andes"; S = 1000; T = 15 },
    @{ Row = 91; D = 44923; K = "Dina"; L = "Segunda";  M = 220; N = 12000; O = 12000; P = 12000; Q = "`$/caja 15 kilos granel"; R = "Provincia de Los Andes"; S =  800; T = 15 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $A
    $ws.Cells.Item($row, 2).Value = $B
    $ws.Cells.Item($row, 3).Value = $C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $E
    $ws.Cells.Item($row, 6).Value = $F
    $ws.Cells.Item($row, 7).Value = $G
    $ws.Cells.Item($row, 8).Value = $H
    $ws.Cells.Item($row, 9).Value = $I
    $ws.Cells.Item($row, 10).Value = $J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
